# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables "magic row" headers on each worksheet encode their
# attributes with UpperCamelCase names (Type=, Id=, Name=, Description=,
# Date=, ObjTablesVersion=, TableID=, TableName=). This change renames
# those attribute keys to lowerCamelCase (type=, id=, name=, description=,
# date=, objTablesVersion=, tableID=, tableName=) while leaving the
# attribute values untouched.

$wb = $excel.ActiveWorkbook

# --- "!!_Table of contents" sheet ---
$tocSheet = $wb.Worksheets.Item("!!_Table of contents")

$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"

$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents' tableID='Table of contents' tableName='Readme' description='Table/model and column/attribute definitions' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# --- "!!Child" sheet ---
$childSheet = $wb.Worksheets.Item("!!Child")

$childSheet.Range("A1").Value = "!!ObjTables type='Data' id='Child' name='Child' description='Represents a child' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# --- "!!Parent" sheet ---
$parentSheet = $wb.Worksheets.Item("!!Parent")

$parentSheet.Range("A1").Value = "!!ObjTables type='Data' id='Parent' name='Parent' description='Represents a parent' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"
